$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Repeats" value in the Values row (I4) from 1 back to 200
$ws.Range("I4").Value = "200"

# Restore the active selection to J13 (matches recorded cursor position)
$ws.Range("J13").Select()
